$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), columns AF..AL ---
$ws.Range("AF1").Value = "timestart"
$ws.Range("AG1").Value = "timeend"
$ws.Range("AH1").Value = "NRO_DIAS_DE_MATRICULAS"
$ws.Range("AI1").Value = "El tiempo de matricula es invalido"
$ws.Range("AJ1").Value = "Numero_Wapp_Incorrecto"
$ws.Range("AK1").Value = "Numero_Con_Prefijo"
$ws.Range("AL1").Value = "El campo del pais esta vacío"

# Copy the header style (bold/border/center) from an existing header cell (AE1)
# onto the newly added header cells so they match the rest of row 1.
$ws.Range("AE1").Copy()
$ws.Range("AF1:AL1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 2 ---
$ws.Range("G2").ClearContents()
$ws.Range("AE2").Value = 25
$ws.Range("AF2").Value = 1725526800
$ws.Range("AG2").Value = 1727686800
$ws.Range("AH2").Value = 25
$ws.Range("AI2").Value = "NO"
$ws.Range("AJ2").Value = "NO"
$ws.Range("AK2").ClearContents()
$ws.Range("AL2").Value = "NO"

# --- Row 3 ---
$ws.Range("AF3").Value = 1725526800
$ws.Range("AG3").Value = 1725526800
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = "NO"
$ws.Range("AJ3").Value = "SI"
$ws.Range("AK3").Value = "+5732090wr0{6290a"
$ws.Range("AL3").Value = "SI"

# --- Row 4 ---
$ws.Range("AE4").Value = 7
$ws.Range("AF4").Value = 1725526800
$ws.Range("AG4").Value = 1726131600
$ws.Range("AH4").Value = 7
$ws.Range("AI4").Value = "NO"
$ws.Range("AJ4").Value = "NO"
$ws.Range("AK4").NumberFormat = "@"
$ws.Range("AK4").Value = "+573209006290"
$ws.Range("AL4").Value = "NO"

# --- Row 5 ---
$ws.Range("AE5").Value = 7
$ws.Range("AF5").Value = 1725526800
$ws.Range("AG5").Value = 1726131600
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = "NO"
$ws.Range("AJ5").Value = "NO"
$ws.Range("AK5").NumberFormat = "@"
$ws.Range("AK5").Value = "+573209006290"
$ws.Range("AL5").Value = "NO"

# --- Row 6 ---
$ws.Range("AE6").Value = 25
$ws.Range("AF6").Value = 1725526800
$ws.Range("AG6").Value = 1726131600
$ws.Range("AH6").Value = 7
$ws.Range("AI6").Value = "NO"
$ws.Range("AJ6").Value = "NO"
$ws.Range("AK6").NumberFormat = "@"
$ws.Range("AK6").Value = "+573183812254"
$ws.Range("AL6").Value = "NO"

# --- Row 7 ---
$ws.Range("AE7").Value = 25
$ws.Range("AF7").Value = 1725526800
$ws.Range("AG7").Value = 1727686800
$ws.Range("AH7").Value = 25
$ws.Range("AI7").Value = "NO"
$ws.Range("AJ7").Value = "NO"
$ws.Range("AK7").NumberFormat = "@"
$ws.Range("AK7").Value = "+573183812254"
$ws.Range("AL7").Value = "NO"
